# Edit script for GO5050_Silvotti.docx
#
# The diff:
#  1. First paragraph (title/author block): pPr drops widowControl /
#     autoSpaceDE / autoSpaceDN / adjustRightInd / spacing and the rPr's
#     bCs, and gains eastAsia="Times New Roman" on rFonts. The
#     "Roberto " + "Silvot" + "ti" runs are re-split into "R" / "oberto "
#     / "Silvotti" with a _GoBack bookmark inserted right after the "R",
#     then a manual line break and a new "INAF" run (affiliation) are
#     appended.
#  2. Second paragraph ("Although only 2% ..."): a manual line break run
#     is prepended before the existing text.
#  3. The last (empty) paragraph loses the _GoBack bookmark (it moved to
#     paragraph 1, per above).
#
# This runtime's Range.InsertXML replaces the *entire* paragraph(s)
# touched by the range with the pasted OOXML body content, so each
# paragraph below is given in full (only the relevant bits differ from
# the source).

$d = $word.ActiveDocument

# --- 1) First paragraph: title/author block -----------------------------
$para1Xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00351E0D" w:rsidRPr="00301FEF" w:rsidRDefault="00301FEF" w:rsidP="00301FEF"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">A survey to detect first </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>sdB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> Planetary Transits </w:t></w:r><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>R</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">oberto </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Silvotti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>INAF</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p1 = $d.Paragraphs(1)
[void]$p1.Range.InsertXML($para1Xml)

# --- 2) Second paragraph: prepend a manual line break --------------------
$para2Xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00301FEF" w:rsidRPr="00301FEF" w:rsidRDefault="00301FEF" w:rsidP="00301FEF"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:br/></w:r><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">Although only 2% of stars evolve through the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>subdwarf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> B (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>sdB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">) phase (e.g. Heber 2009 ARA&amp;A 47, 211 and refs. therein), there are at least two good reasons to study </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>sdB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00301FEF"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> planets/BDs: </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p2 = $d.Paragraphs(2)
[void]$p2.Range.InsertXML($para2Xml)

# --- 3) Last (empty) paragraph: drop the old _GoBack bookmark -----------
# Replacing the very last paragraph in the story leaves the document's
# sacrosanct final paragraph mark behind, which would add a spurious
# extra empty paragraph. Insert the fixed-up paragraph, then delete the
# now-redundant duplicate it pushed ahead of the real last mark.
$lastParaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00301FEF" w:rsidRPr="00301FEF" w:rsidRDefault="00301FEF" w:rsidP="00301FEF"><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$countBefore = $d.Paragraphs.Count
$pLast = $d.Paragraphs($countBefore)
[void]$pLast.Range.InsertXML($lastParaXml)
if ($d.Paragraphs.Count -gt $countBefore) {
    $d.Paragraphs($countBefore).Range.Delete()
}

Write-Output "done"
